# Fill in the empty attendance cells (columns A-H, rows 3-35) with "-"
# so that the attendance sheet accurately reflects absences ("-") versus
# attendances ("+"), matching the accurate broj_izostanka (absence count).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:C3").Value = "-"
$ws.Range("G3").Value = "-"
$ws.Range("G4").Value = "-"
$ws.Range("B5").Value = "-"
$ws.Range("E5:G5").Value = "-"
$ws.Range("B6:C6").Value = "-"
$ws.Range("B7").Value = "-"
$ws.Range("B8:H8").Value = "-"
$ws.Range("C10").Value = "-"
$ws.Range("G10:H10").Value = "-"
$ws.Range("B12:H12").Value = "-"
$ws.Range("B13:C13").Value = "-"
$ws.Range("E13:G13").Value = "-"
$ws.Range("B14:H14").Value = "-"
$ws.Range("B16").Value = "-"
$ws.Range("E16:G16").Value = "-"
$ws.Range("B17").Value = "-"
$ws.Range("B18:C18").Value = "-"
$ws.Range("F18:G18").Value = "-"
$ws.Range("B19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("H19").Value = "-"
$ws.Range("B20:C20").Value = "-"
$ws.Range("E20:H20").Value = "-"
$ws.Range("B21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("H21").Value = "-"
$ws.Range("D22").Value = "-"
$ws.Range("H22").Value = "-"
$ws.Range("G23:H23").Value = "-"
$ws.Range("B24").Value = "-"
$ws.Range("G25").Value = "-"
$ws.Range("D26").Value = "-"
$ws.Range("B27:C27").Value = "-"
$ws.Range("E27").Value = "-"
$ws.Range("G27").Value = "-"
$ws.Range("D28").Value = "-"
$ws.Range("F28:H28").Value = "-"
$ws.Range("B29:C29").Value = "-"
$ws.Range("G29").Value = "-"
$ws.Range("D30").Value = "-"
$ws.Range("B31:E31").Value = "-"
$ws.Range("G31").Value = "-"
$ws.Range("A32:H32").Value = "-"
$ws.Range("A33:H33").Value = "-"
$ws.Range("A34:H34").Value = "-"
$ws.Range("A35:H35").Value = "-"
